$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 11.14494766666667
$ws.Range("H2").Value = 33.434843
$ws.Range("I2").Value = 0.1279818847384872
$ws.Range("J2").Value = 0.1279818847384872
$ws.Range("Q2").Value = 10.29505253251944
$ws.Range("R2").Value = 92.655472792675
$ws.Range("S2").Value = 0.1251271758759513
$ws.Range("T2").Value = 0.1251271758759513

# Row 3
$ws.Range("G3").Value = 11.14494766666667
$ws.Range("H3").Value = 33.434843
$ws.Range("I3").Value = 0.1279818847384872
$ws.Range("J3").Value = 0.1279818847384872
$ws.Range("Q3").Value = 0.2348760570924445
$ws.Range("R3").Value = 2.113884513832
$ws.Range("S3").Value = 0.002854708862535935
$ws.Range("T3").Value = 0.002854708862535935

# Row 4
$ws.Range("I4").Value = 0.5307607770439682
$ws.Range("J4").Value = 0.5307607770439681
$ws.Range("S4").Value = 0.5189218554871408
$ws.Range("T4").Value = 0.5189218554871406

# Row 5
$ws.Range("I5").Value = 0.5307607770439682
$ws.Range("J5").Value = 0.5307607770439681
$ws.Range("S5").Value = 0.01183892155682739
$ws.Range("T5").Value = 0.01183892155682739

# Row 6
$ws.Range("G6").Value = 29.71744933333333
$ws.Range("H6").Value = 89.152348
$ws.Range("I6").Value = 0.3412573382175446
$ws.Range("J6").Value = 0.3412573382175446
$ws.Range("Q6").Value = 27.45124617625556
$ws.Range("R6").Value = 247.0612155863
$ws.Range("S6").Value = 0.3336453988418613
$ws.Range("T6").Value = 0.3336453988418613

# Row 7
$ws.Range("G7").Value = 29.71744933333333
$ws.Range("H7").Value = 89.152348
$ws.Range("I7").Value = 0.3412573382175446
$ws.Range("J7").Value = 0.3412573382175446
$ws.Range("Q7").Value = 0.6262853388835556
$ws.Range("R7").Value = 5.636568049952
$ws.Range("S7").Value = 0.007611939375683261
$ws.Range("T7").Value = 0.007611939375683261
